$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 76, shifting existing rows (76..189) down to (77..190).
$ws.Rows(76).Insert()

# Populate the newly inserted row 76 with a new price record.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R carry the same "template" values
# as the surrounding records for this market/category; D, J, K, L, M, P hold
# the new observation's data.
$ws.Cells.Item(76, 1).Value = 3
$ws.Cells.Item(76, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(76, 3).Value = "Coquimbo"
$ws.Cells.Item(76, 4).Value = 44803
$ws.Cells.Item(76, 5).Value = 5
$ws.Cells.Item(76, 6).Value = 100112026
$ws.Cells.Item(76, 7).Value = "Haba"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 105
$ws.Cells.Item(76, 11).Value = 13000
$ws.Cells.Item(76, 12).Value = 14000
$ws.Cells.Item(76, 13).Value = 13524
$ws.Cells.Item(76, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(76, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(76, 16).Value = 541
$ws.Cells.Item(76, 17).Value = 25
$ws.Cells.Item(76, 18).Value = "Hortaliza"
